# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 427 of the data sheet,
# pushing all existing rows (old 427-456) down by one (new 428-457).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 427 - shifts rows 427:456 down to 428:457
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with the new record
$ws.Range("A427").Value = 11
$ws.Range("B427").Value = "Vega Monumental Concepción"
$ws.Range("C427").Value = "Bíobío"
$ws.Range("D427").Value = 45021
$ws.Range("E427").Value = 8
$ws.Range("F427").Value = "Fruta"
$ws.Range("G427").Value = 100102
$ws.Range("H427").Value = "Cítricos"
$ws.Range("I427").Value = 100102005
$ws.Range("J427").Value = "Naranja"
$ws.Range("K427").Value = "Valencia"
$ws.Range("L427").Value = "Primera"
$ws.Range("M427").Value = 400
$ws.Range("N427").Value = 13000
$ws.Range("O427").Value = 14000
$ws.Range("P427").Value = 13625
$ws.Range("Q427").Value = "`$/caja 15 kilos granel"
$ws.Range("R427").Value = "Región de O'Higgins"
$ws.Range("S427").Value = 908
$ws.Range("T427").Value = 15
